$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I: "MVP 2.0" header, marked "x" for the 5 task rows, and filled-by note at row 8
$ws.Range("I1").Value = "MVP 2.0"
$ws.Range("I2:I6").Value = "x"

# Copy the style of H8 ("Ausgefüllt von" row) onto I8, then set its value
$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I8").Value = "Pütter/Hesse"

[void]$ws.Range("I9").Select()
